$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style of D1 to E1 (bold/centered/bordered header style) then set header value
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Cells.Item(1, 5).Value = 3

# Column E values for rows 2-231 (Vietnamese translation mirroring column C grouping)
$colE = @{
    2 = "lưu chuyển tiền từ hoạt động kinh doanh"
    3 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    4 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    5 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    6 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    7 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    8 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    9 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    10 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    11 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    12 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    13 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    14 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    15 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    16 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    17 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    18 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    19 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    20 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    21 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    22 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    23 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    24 = "lưu chuyển tiền thuần trong kỳ"
    25 = "lưu chuyển tiền từ hoạt động đầu tư"
    26 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    27 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    28 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    29 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    30 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    31 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    32 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    33 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    34 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    35 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    36 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    37 = "lưu chuyển tiền thuần trong kỳ"
    38 = "lưu chuyển tiền từ hoạt động tài chính"
    39 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    40 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    41 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    42 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    43 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    44 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    45 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    46 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    47 = "lưu chuyển tiền thuần trong kỳ"
    48 = "tiền và tương đương tiền cuối kỳ"
    49 = "tiền và tương đương tiền cuối kỳ"
    50 = "tiền và tương đương tiền cuối kỳ"
    51 = "tiền và tương đương tiền cuối kỳ"
    52 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    53 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    54 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    55 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    56 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    57 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    58 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    59 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    60 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    61 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    62 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    63 = "lưu chuyển tiền từ hoạt động kinh doanh"
    64 = "3 lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    65 = "3 lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    66 = "3 lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    67 = "3 lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    68 = "3 lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    69 = "3 lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    70 = "3 lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    71 = "3 lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    72 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    73 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    74 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    75 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    76 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    77 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    78 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    79 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    80 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    81 = "lưu chuyển tiền thuần từ hoạt động kinh doanh"
    82 = "lưu chuyển tiền thuần trong kỳ"
    83 = "lưu chuyển tiền từ hoạt động quản lý danh mục đầu tư cho người đầu tư"
    84 = "lưu chuyển tiền thuần từ hoạt động quản lý danh mục đầu tư cho người đầu tư"
    85 = "lưu chuyển tiền thuần từ hoạt động quản lý danh mục đầu tư cho người đầu tư"
    86 = "lưu chuyển tiền thuần từ hoạt động quản lý danh mục đầu tư cho người đầu tư"
    87 = "lưu chuyển tiền thuần từ hoạt động quản lý danh mục đầu tư cho người đầu tư"
    88 = "lưu chuyển tiền thuần từ hoạt động quản lý danh mục đầu tư cho người đầu tư"
    89 = "lưu chuyển tiền thuần từ hoạt động quản lý danh mục đầu tư cho người đầu tư"
    90 = "lưu chuyển tiền thuần từ hoạt động quản lý danh mục đầu tư cho người đầu tư"
    91 = "lưu chuyển tiền thuần từ hoạt động quản lý danh mục đầu tư cho người đầu tư"
    92 = "lưu chuyển tiền thuần trong kỳ"
    93 = "lưu chuyển tiền từ hoạt động đầu tư"
    94 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    95 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    96 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    97 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    98 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    99 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    100 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    101 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    102 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    103 = "lưu chuyển tiền thuần trong kỳ"
    104 = "lưu chuyển tiền từ hoạt động tài chính"
    105 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    106 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    107 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    108 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    109 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    110 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    111 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    112 = "lưu chuyển tiền thuần trong kỳ"
    113 = "tiền và tương đương tiền cuối kỳ"
    114 = "tiền và tương đương tiền cuối kỳ"
    115 = "tiền và tương đương tiền cuối kỳ"
    116 = "tiền và tương đương tiền cuối kỳ"
    117 = "lưu chuyển tiền từ hoạt động kinh doanh chứng khoán"
    118 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    119 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    120 = "tăng các chi phí phi tiền tệ"
    121 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    122 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    123 = "tăng các chi phí phi tiền tệ"
    124 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    125 = "tăng các chi phí phi tiền tệ"
    126 = "tăng các chi phí phi tiền tệ"
    127 = "tăng các chi phí phi tiền tệ"
    128 = "tăng các chi phí phi tiền tệ"
    129 = "tăng các chi phí phi tiền tệ"
    130 = "tăng các chi phí phi tiền tệ"
    131 = "tăng các chi phí phi tiền tệ"
    132 = "tăng các chi phí phi tiền tệ"
    133 = "tăng các chi phí phi tiền tệ"
    134 = "tăng các chi phí phi tiền tệ"
    135 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    136 = "giảm các doanh thu phi tiền tệ"
    137 = "giảm các doanh thu phi tiền tệ"
    138 = "giảm các doanh thu phi tiền tệ"
    139 = "giảm các doanh thu phi tiền tệ"
    140 = "giảm các doanh thu phi tiền tệ"
    141 = "giảm các doanh thu phi tiền tệ"
    142 = "giảm các doanh thu phi tiền tệ"
    143 = "giảm các doanh thu phi tiền tệ"
    144 = "giảm các doanh thu phi tiền tệ"
    145 = "giảm các doanh thu phi tiền tệ"
    146 = "lợi nhuận từ hoạt động kinh doanh trước thay đổi vốn lưu động"
    147 = "tăng các chi phí phi tiền tệ"
    148 = "tăng các chi phí phi tiền tệ"
    149 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    150 = "tăng các chi phí phi tiền tệ"
    151 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    152 = "tăng các chi phí phi tiền tệ"
    153 = "tăng các chi phí phi tiền tệ"
    154 = "tăng các chi phí phi tiền tệ"
    155 = "tăng các chi phí phi tiền tệ"
    156 = "tăng các chi phí phi tiền tệ"
    157 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    158 = "tăng các chi phí phi tiền tệ"
    159 = "tăng các chi phí phi tiền tệ"
    160 = "tăng các chi phí phi tiền tệ"
    161 = "tăng các chi phí phi tiền tệ"
    162 = "tăng các chi phí phi tiền tệ"
    163 = "tăng các chi phí phi tiền tệ"
    164 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    165 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    166 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    167 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    168 = "tăng các chi phí phi tiền tệ"
    169 = "tăng các chi phí phi tiền tệ"
    170 = "tăng giảm các khoản trích nộp phúc lợi nhân viên"
    171 = "tăng các chi phí phi tiền tệ"
    172 = "tăng giảm các khoản trích nộp phúc lợi nhân viên"
    173 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    174 = "tăng các chi phí phi tiền tệ"
    175 = "tăng các chi phí phi tiền tệ"
    176 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    177 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    178 = "lưu chuyển tiền thuần từ hoạt động kinh doanh chứng khoán"
    179 = "lưu chuyển tiền thuần trong kỳ"
    180 = "lưu chuyển tiền từ hoạt động đầu tư"
    181 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    182 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    183 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    184 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    185 = "lưu chuyển tiền thuần từ hoạt động đầu tư"
    186 = "lưu chuyển tiền từ hoạt động tài chính"
    187 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    188 = "lưu chuyển tiền thuần từ hoạt động tài chính"
    189 = "tăng các chi phí phi tiền tệ"
    190 = "tăng các chi phí phi tiền tệ"
    191 = "tăng các chi phí phi tiền tệ"
    192 = "tăng các chi phí phi tiền tệ"
    193 = "tăng các chi phí phi tiền tệ"
    194 = "tiền và tương đương tiền cuối kỳ"
    195 = "tăng các chi phí phi tiền tệ"
    196 = "tăng các chi phí phi tiền tệ"
    197 = "tăng các chi phí phi tiền tệ"
    198 = "tăng các chi phí phi tiền tệ"
    199 = "tăng các chi phí phi tiền tệ"
    200 = "tăng các chi phí phi tiền tệ"
    201 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    202 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    203 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    204 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    205 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    206 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    207 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    208 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    209 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    210 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    211 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    212 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    213 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    214 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    215 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    216 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    217 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    218 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    219 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    220 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    221 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    222 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    223 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    224 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    225 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    226 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    227 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    228 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    229 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    230 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
    231 = "phần lưu chuyển tiền tệ hoạt động môi giới ủy thác của khách hàng"
}

foreach ($r in $colE.Keys) {
    $ws.Cells.Item($r, 5).Value = $colE[$r]
}

